# Refresh cryptocurrency price/volume snapshot (GitHub Actions data pull).
# Row order for three coin pairs also rotated (new rank order), so Coin/Link
# cells are rewritten for those rows too.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ col = new value }  (only cells that actually changed)
$updates = @{
    2 = @{ 'D'='27.069.61'; 'E'='  -0.11%  ' }
    3 = @{ 'D'='1.884.83'; 'E'='  -0.53%  ' }
    4 = @{ 'D'='1.001'; 'E'='  -0.27%  ' }
    5 = @{ 'D'='304.27'; 'E'='  -0.99%  ' }
    6 = @{ 'D'='1.001'; 'E'='  -0.20%  ' }
    7 = @{ 'D'='0.5375'; 'E'='  +3.57%  ' }
    8 = @{ 'D'='0.3780'; 'E'='  +0.17%  ' }
    9 = @{ 'D'='0.07240'; 'E'='  -0.57%  ' }
    10 = @{ 'D'='22.02'; 'E'='  +4.06%  ' }
    11 = @{ 'D'='0.8911'; 'E'='  -0.88%  ' }
    12 = @{ 'D'='0.08131'; 'E'='  -0.52%  ' }
    13 = @{ 'B'='WrappedEther'; 'C'='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; 'D'='1.893.60'; 'E'='  -0.39%  ' }
    14 = @{ 'B'='Litecoin'; 'C'='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; 'D'='93.90'; 'E'='  -1.55%  ' }
    15 = @{ 'B'='Polkadot'; 'C'='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; 'D'='5.306'; 'E'='  -0.64%  ' }
    16 = @{ 'D'='1.001'; 'E'='  -0.32%  ' }
    17 = @{ 'D'='14.81'; 'E'='  +1.46%  ' }
    18 = @{ 'D'='0.000008607'; 'E'='  +0.13%  ' }
    19 = @{ 'D'='1.001'; 'E'='  -0.27%  ' }
    20 = @{ 'D'='27.072.57'; 'E'='  -0.19%  ' }
    21 = @{ 'D'='5.005'; 'E'='  -1.74%  ' }
    22 = @{ 'D'='10.79'; 'E'='  +0.78%  ' }
    23 = @{ 'D'='6.457'; 'E'='  +0.29%  ' }
    24 = @{ 'D'='148.28'; 'E'='  -0.63%  ' }
    25 = @{ 'D'='2.291'; 'E'='  -1.14%  ' }
    26 = @{ 'B'='Toncoin'; 'C'='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; 'D'='1.747'; 'E'='  +0.25%  ' }
    27 = @{ 'B'='EthereumClassic'; 'C'='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; 'D'='18.18'; 'E'='  -0.21%  ' }
    28 = @{ 'D'='115.58' }
    29 = @{ 'D'='4.766'; 'E'='  -0.58%  ' }
    30 = @{ 'D'='4.637' }
    31 = @{ 'D'='0.09168'; 'E'='  -0.34%  ' }
    32 = @{ 'D'='0.8105'; 'E'='  +2.28%  ' }
    33 = @{ 'D'='0.05043'; 'E'='  +0.26%  ' }
    34 = @{ 'D'='1.185'; 'E'='  -2.01%  ' }
    35 = @{ 'D'='3.006'; 'E'='  +1.02%  ' }
    36 = @{ 'B'='TheSandbox'; 'C'='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; 'D'='0.5981'; 'E'='  +5.33%  ' }
    37 = @{ 'D'='2.641'; 'E'='  +1.72%  ' }
    38 = @{ 'B'='MXToken'; 'C'='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; 'D'='3.232'; 'E'='  -5.18%  ' }
    39 = @{ 'D'='0.01976'; 'E'='  -0.78%  ' }
    40 = @{ 'D'='1.073' }
    41 = @{ 'D'='9.014'; 'E'='  +1.20%  ' }
    42 = @{ 'D'='6.557'; 'E'='  -0.49%  ' }
    43 = @{ 'D'='115.54'; 'E'='  -0.95%  ' }
    44 = @{ 'D'='0.5089'; 'E'='  +5.05%  ' }
    45 = @{ 'D'='0.1513'; 'E'='  +0.28%  ' }
    46 = @{ 'D'='1.001'; 'E'='  -0.34%  ' }
    47 = @{ 'D'='10.11'; 'E'='  +0.65%  ' }
    48 = @{ 'D'='1.616'; 'E'='  +0.02%  ' }
    49 = @{ 'D'='37.66'; 'E'='  -1.49%  ' }
    50 = @{ 'D'='0.06058'; 'E'='  +2.05%  ' }
    51 = @{ 'D'='62.74'; 'E'='  -1.83%  ' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cellRef = "$col$row"
        $value = $updates[$row][$col]
        if ($col -eq "D") {
            # Price column holds dotted numeric-looking strings (e.g. "1.893.60",
            # "1.001") that Excel would otherwise coerce to a Number. Force Text
            # format for the write, then drop back to the default style so the
            # cell keeps no explicit style (matches the rest of the sheet).
            $ws.Range($cellRef).NumberFormat = "@"
            $ws.Range($cellRef).Value = $value
            $ws.Range($cellRef).Style = "Normal"
        } else {
            $ws.Range($cellRef).Value = $value
        }
    }
}
